$d = $word.ActiveDocument

# --- RMSE cell: "1799.6028" -> "1799." + "9457"  (split into two runs) ---
$rng = $d.Content
$rng.Find.Execute("1799.6028")
$start = $rng.Start

$r1 = $d.Range($start, $start + 5)
$r1.Text = "1799."

$r2 = $d.Range($start + 5, $start + 9)
$r2.Font.Bold = 1
$r2.Text = "9457"
$r2 = $d.Range($start + 5, $start + 9)
$r2.Font.Bold = 0

# --- MAE cell: "1362.5385" -> "136" + "4" + "." + "072"  (split into four runs) ---
$rng2 = $d.Content
$rng2.Find.Execute("1362.5385")
$start2 = $rng2.Start

$rA = $d.Range($start2, $start2 + 3)
$rA.Text = "136"

$rB = $d.Range($start2 + 3, $start2 + 4)
$rB.Font.Bold = 1
$rB.Text = "4"
$rB = $d.Range($start2 + 3, $start2 + 4)
$rB.Font.Bold = 0

$rC = $d.Range($start2 + 4, $start2 + 5)
$rC.Font.Bold = 1
$rC.Text = "."
$rC = $d.Range($start2 + 4, $start2 + 5)
$rC.Font.Bold = 0

$rD = $d.Range($start2 + 5, $start2 + 9)
$rD.Font.Bold = 1
$rD.Text = "072"
$rD = $d.Range($start2 + 5, $start2 + 8)
$rD.Font.Bold = 0
